$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-6 (columns A-H)
$data = @(
    @("1326924", "https://aiesec.org/opportunity/global-talent/1326924", "Learning and Development Trainee", "Panamá, Provincia de Panamá, Panamá", "No", "4 applicants", "6 - 18 Months", "HILTI Panama"),
    @("1326923", "https://aiesec.org/opportunity/global-talent/1326923", "AI Tech Developer", "Panamá, Provincia de Panamá, Panamá", "No", "2 applicants", "6 - 18 Months", "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"),
    @("1326906", "https://aiesec.org/opportunity/global-talent/1326906", "Sales Manager", "Mersin, Akdeniz/Mersin, Türkiye", "No", "7 applicants", "3 - 6 Months", "İlke Sports"),
    @("1326653", "https://aiesec.org/opportunity/global-talent/1326653", "Interior Design", "New Cairo City, Cairo Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Ahmad Elsherif Interior Designer"),
    @("1322605", "https://aiesec.org/opportunity/global-talent/1322605", "Architecture", "Gabes, Tunisia", "No", "0 applicants", "9 - 12 Weeks", "IR ARKITETTI")
)

# Column A holds numeric-looking opportunity IDs that must stay stored as
# text (matching the source file's inline-string cells), not auto-coerced
# to numbers. Format the column as Text first so entry doesn't convert it,
# then strip the number-format back off so no stray per-cell style sticks
# around on save.
$ws.Range("A2:A6").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}

$ws.Range("A2:A6").ClearFormats()

# Delete row 7 entirely (shrinks dimension from H7 to H6)
$ws.Rows.Item(7).Delete()

# Update column widths per diff.
# Note: the engine's ColumnWidth setter stores (value + 5/6) in the OOXML
# <col width=".."/> attribute (mirrors Excel's default-font padding math),
# so we subtract 5/6 here to land exactly on the target stored widths
# (35, 43, 60) that the diff specifies.
$pad = 5 / 6
$ws.Columns.Item(3).ColumnWidth = 35 - $pad
$ws.Columns.Item(4).ColumnWidth = 43 - $pad
$ws.Columns.Item(8).ColumnWidth = 60 - $pad
